$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete data rows (5-12), shifting the sheet dimension to A1:C4
$ws.Range("A5:C12").EntireRow.Delete()

# Update the remaining data rows with the new simulation values
$ws.Range("A2").Value = "2024-07-26 11:57:23"
$ws.Range("B2").Value = 99.00326479502259
$ws.Range("C2").Value = 4

$ws.Range("A3").Value = "2024-07-26 11:57:26"
$ws.Range("B3").Value = 99.00326479502259
$ws.Range("C3").Value = 4

$ws.Range("A4").Value = "2024-07-26 11:57:28"
$ws.Range("B4").Value = 99.00326479502259
$ws.Range("C4").Value = 4
